$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data for "Luyện tập mảng 02" exercise, rows 65-70
$labels = @(
    "Luyện tập mảng 02 - 01",
    "Luyện tập mảng 02 - 02",
    "Luyện tập mảng 02 - 03",
    "Luyện tập mảng 02 - 04",
    "Luyện tập mảng 02 - 05",
    "Luyện tập mảng 02 - 06"
)

$urls = @(
    "https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s12-mang/exercise/%5BB%C3%A0i%20t%E1%BA%ADp%201%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20m%E1%BA%A3ng%2002.html",
    "https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s12-mang/exercise/%5BB%C3%A0i%20t%E1%BA%ADp%202%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20m%E1%BA%A3ng%2002.html",
    "https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s12-mang/exercise/%5BB%C3%A0i%20t%E1%BA%ADp%203%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20m%E1%BA%A3ng%2002.html",
    "https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s12-mang/exercise/%5BB%C3%A0i%20t%E1%BA%ADp%204%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20m%E1%BA%A3ng%2002.html",
    "https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s12-mang/exercise/%5BB%C3%A0i%20t%E1%BA%ADp%205%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20m%E1%BA%A3ng%2002.html",
    "https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s12-mang/exercise/%5BB%C3%A0i%20t%E1%BA%ADp%206%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20m%E1%BA%A3ng%2002.html"
)

$startRow = 65

# Replicate the original authoring order so the shared-strings table is
# built in the same sequence as the source workbook:
#   1) B65 (first url)
#   2) A65:A70 (all labels)
#   3) B66:B70 (remaining urls)
$ws.Cells.Item($startRow, 2).Value = $urls[0]

for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $labels[$i]
}

for ($i = 1; $i -lt $urls.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $urls[$i]
}

# Update the visible window / selection to match the new bottom of the sheet
$excel.ActiveWindow.ScrollRow = 54
$ws.Range("B65:B70").Select()

$wb.Save()
